$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.365.88'
$ws.Range("E2").Value = '  +0.28%  '

$ws.Range("D3").Value = '1.869.92'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7110'
$ws.Range("E5").Value = '  +1.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '238.34'
$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07795'
$ws.Range("E8").Value = '  -5.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3064'
$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.11'
$ws.Range("E10").Value = '  +7.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08181'
$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("D12").Value = '1.881.72'
$ws.Range("E12").Value = '  +0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.246'
$ws.Range("E13").Value = '  +1.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7196'
$ws.Range("E14").Value = '  +0.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.25'
$ws.Range("E15").Value = '  -0.06%  '

$ws.Range("D16").Value = '29.364.19'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.827'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '241.61'
$ws.Range("E18").Value = '  +1.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007816'
$ws.Range("E19").Value = '  -0.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.26'
$ws.Range("E20").Value = '  -0.94%  '

$ws.Range("D21").Value = '2.125.91'
$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.621'
$ws.Range("E24").Value = '  +2.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '162.32'
$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.957'
$ws.Range("E26").Value = '  -0.32%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1452'
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.18'
$ws.Range("E28").Value = '  +0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.921'
$ws.Range("E29").Value = '  -3.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.360'
$ws.Range("E30").Value = '  -5.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.477'
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.330'
$ws.Range("E32").Value = '  -2.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.056'
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05220'
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.189'
$ws.Range("E35").Value = '  +1.48%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7173'
$ws.Range("E36").Value = '  +1.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.675'
$ws.Range("E38").Value = '  +0.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01852'
$ws.Range("E39").Value = '  -0.03%  '

$ws.Range("E40").Value = '  -0.83%  '

$ws.Range("D41").Value = '1.180.30'
$ws.Range("E41").Value = '  +2.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9159'
$ws.Range("E42").Value = '  -0.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.007'
$ws.Range("E43").Value = '  +0.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4293'
$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.42'
$ws.Range("E45").Value = '  +0.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.17'
$ws.Range("E47").Value = '  -0.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5362'
$ws.Range("E48").Value = '  -0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.764'
$ws.Range("E49").Value = '  -0.50%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.251'
$ws.Range("E50").Value = '  +0.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.041'
$ws.Range("E51").Value = '  +0.95%  '
